$p = $ppt.ActivePresentation

# --- Update "Date Placeholder" fields on the slide master and every slide layout ---
# (source deck was re-saved two days later: 23/05/2023 -> 25/05/2023)
function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "25/05/2023"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes
$layouts = $master.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DatePlaceholders $layouts.Item($j).Shapes
}

# --- Slide 2 (sldId 257, cId 841685169): shape 3 "Marcador de Posicao de Conteudo 2" ---
# Replace the short "Objetivos" placeholder text with the full project description paragraph.
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(2)
$tr2 = $sh2.TextFrame.TextRange
$tr2.Text = 'O '
$cur2 = $tr2
$cur2 = $cur2.InsertAfter('projeto')
$cur2 = $cur2.InsertAfter(' da ')
$cur2 = $cur2.InsertAfter('Unidade')
$cur2 = $cur2.InsertAfter(' Curricular ')
$cur2 = $cur2.InsertAfter('Metodologias')
$cur2 = $cur2.InsertAfter(' de ')
$cur2 = $cur2.InsertAfter('Desenvilvimento')
$cur2 = $cur2.InsertAfter(' de Software ')
$cur2 = $cur2.InsertAfter('consiste')
$cur2 = $cur2.InsertAfter(' ')
$cur2 = $cur2.InsertAfter('na')
$cur2 = $cur2.InsertAfter(' ')
$cur2 = $cur2.InsertAfter('gestão')
$cur2 = $cur2.InsertAfter(' dos ')
$cur2 = $cur2.InsertAfter('projetos')
$cur2 = $cur2.InsertAfter(' das UC’s de ')
$cur2 = $cur2.InsertAfter('Desenvolvimento')
$cur2 = $cur2.InsertAfter(' de ')
$cur2 = $cur2.InsertAfter('Aplicações')
$cur2 = $cur2.InsertAfter(' e de ')
$cur2 = $cur2.InsertAfter('Programação')
$cur2 = $cur2.InsertAfter(' para a WEB – ')
$cur2 = $cur2.InsertAfter('Servidor')
$cur2 = $cur2.InsertAfter(' ')
$cur2 = $cur2.InsertAfter('através')
$cur2 = $cur2.InsertAfter(' de ')
$cur2 = $cur2.InsertAfter('Metodologias')
$cur2 = $cur2.InsertAfter(' de ')
$cur2 = $cur2.InsertAfter('Desenvolvimento')
$cur2 = $cur2.InsertAfter(' ')
$cur2 = $cur2.InsertAfter('Ágies')
$cur2 = $cur2.InsertAfter(' e a ')
$cur2 = $cur2.InsertAfter('Metodologia')
$cur2 = $cur2.InsertAfter(' de SCRUM.')
$cur2 = $cur2.InsertAfter("`r")

# --- Slide 4 (sldId 259, cId 667560872): shape 8 "Marcador de Posicao de Conteudo 2" ---
# Fix the acronym typo "PWCS" -> "PWS".
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(2)
$tr4 = $sh4.TextFrame.TextRange
$full4 = $tr4.Text
$idx4 = $full4.IndexOf("  de PWCS ")
$run4 = $tr4.Characters($idx4 + 1, 10)
$run4.Text = '  de PWS '

